$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "RM 232" row (original row 26) and the "SC 92" row (original row 28)
# were removed from the data set. Deleting higher row number first keeps
# the lower row index valid.
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Cells whose numbers were newly imputed (previously blank/missing)
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("E15").Value = -8.4
$ws.Range("E25").Value = -7.1

# Cells that became missing (previously had a value, now blank)
$ws.Range("C3").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("E19").ClearContents()
$ws.Range("E29").ClearContents()

# Value corrections on the rows shifted up after the deletions above
# (new row 26 = "SC 5", new row 27 = "SC 101", new row 33 = "SC 232")
$ws.Range("B26").Value = -20.2
$ws.Range("B27").ClearContents()
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
